$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Refresh the "Price" (D) column. Most prices look like plain decimal
# numbers to Excel's type inference, so force the cell to Text first
# (then drop the style back to Normal so no visible formatting changes),
# otherwise e.g. "1.00" would be auto-coerced to the number 1.
$ws.Range("D2").Value = "61.589.87"
$ws.Range("D3").Value = "2.893.16"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "572.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.16"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.503"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Value = "2.895.79"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.69"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.149"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000233"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.25"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Value = "3.373.10"
$ws.Range("D17").Value = "61.667.24"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Value = "2.890.32"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.75"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.655"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.90"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.53"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "11.82"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.03"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000107"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.01"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.52"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.07"
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.106"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "25.54"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.959"
$ws.Range("D36").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.95"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.115"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "8.26"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.268"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "38.36"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Value = "2.677.17"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "133.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0333"
$ws.Range("D47").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "21.66"
$ws.Range("D51").Style = "Normal"

# Refresh the "Volume(1h)" (E) column. These are already non-numeric
# text (padded percentages), so a plain .Value assignment is safe.
$ws.Range("E2").Value = "  -2.03%  "
$ws.Range("E3").Value = "  -2.49%  "
$ws.Range("E4").Value = "  +0.31%  "
$ws.Range("E5").Value = "  -3.88%  "
$ws.Range("E6").Value = "  -1.23%  "
$ws.Range("E7").Value = "  +0.14%  "
$ws.Range("E8").Value = "  -0.44%  "
$ws.Range("E9").Value = "  -2.29%  "
$ws.Range("E10").Value = "  -8.40%  "
$ws.Range("E11").Value = "  -0.49%  "
$ws.Range("E12").Value = "  -3.26%  "
$ws.Range("E13").Value = "  -2.95%  "
$ws.Range("E14").Value = "  -2.66%  "
$ws.Range("E15").Value = "  -0.71%  "
$ws.Range("E16").Value = "  -2.43%  "
$ws.Range("E17").Value = "  -1.65%  "
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("E19").Value = "  -2.67%  "
$ws.Range("E20").Value = "  -2.18%  "
$ws.Range("E21").Value = "  -1.91%  "
$ws.Range("E22").Value = "  -2.56%  "
$ws.Range("E23").Value = "  -2.44%  "
$ws.Range("E24").Value = "  -2.53%  "
$ws.Range("E25").Value = "  -0.44%  "
$ws.Range("E26").Value = "  -10.33%  "
$ws.Range("E27").Value = "  +0.03%  "
$ws.Range("E28").Value = "  -4.76%  "
$ws.Range("E29").Value = "  +11.04%  "
$ws.Range("E30").Value = "  -3.55%  "
$ws.Range("E31").Value = "  -3.60%  "
$ws.Range("E32").Value = "  -3.41%  "
$ws.Range("E33").Value = "  +0.31%  "
$ws.Range("E34").Value = "  -3.74%  "
$ws.Range("E35").Value = "  -3.75%  "
$ws.Range("E36").Value = "  -3.52%  "
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  -4.48%  "
$ws.Range("E41").Value = "  -2.05%  "
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("E43").Value = "  -4.60%  "
$ws.Range("E44").Value = "  -5.07%  "
$ws.Range("E45").Value = "  -2.35%  "
$ws.Range("E46").Value = "  -1.18%  "
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("E50").Value = "  -2.07%  "
$ws.Range("E51").Value = "  -6.04%  "

# Rows 37/38 and 48/49 swapped ranking order, each also with refreshed
# price/volume figures, so update every column for these four rows.

# Row 37: dogwifhat -> Filecoin
$ws.Range("B37").Value = "Filecoin"
$ws.Range("C37").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.43"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.82%  "

# Row 38: Filecoin -> dogwifhat
$ws.Range("B38").Value = "dogwifhat"
$ws.Range("C38").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.00"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.21%  "

# Row 48: Bittensor -> USDe
$ws.Range("B48").Value = "USDe"
$ws.Range("C48").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.00%  "

# Row 49: USDe -> Bittensor
$ws.Range("B49").Value = "Bittensor"
$ws.Range("C49").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "340.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.37%  "
